$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("parametrosInicio")

# Update the migration label text in B10
$ws.Range("B10").Value = "MIGRACIONES SGV MARZO 2023 15.03.2023 BOOT 5"

# Auto-fit column B to the new (longer) text, matching bestFit column width behavior
$ws.Columns.Item(2).EntireColumn.AutoFit() | Out-Null

# Update the active selection on this sheet to D7
$ws.Activate()
$ws.Range("D7").Select() | Out-Null
